$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay as text (avoid Excel auto-converting
# strings like "1.00" or "611.90" into numbers and losing formatting).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.164.95"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.725.85"
$ws.Range("E3").Value = "  -0.40%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "611.90"
$ws.Range("E5").Value = "  +5.00%  "

# Row 6 - Solana
$ws.Range("D6").Value = "189.48"
$ws.Range("E6").Value = "  +5.98%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.74%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.65%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "59.23"
$ws.Range("E11").Value = "  +9.26%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -3.47%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "10.69"
$ws.Range("E13").Value = "  -1.48%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.316.72"
$ws.Range("E14").Value = "  -0.36%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.722.46"
$ws.Range("E15").Value = "  -1.06%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "19.35"
$ws.Range("E16").Value = "  -1.06%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  -0.89%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -1.34%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.992.69"
$ws.Range("E20").Value = "  +0.95%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "412.61"
$ws.Range("E21").Value = "  -0.16%  "

# Row 22 - PancakeSwap
$ws.Range("E22").Value = "  -0.02%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "89.44"
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -1.64%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "12.86"
$ws.Range("E25").Value = "  -1.12%  "

# Row 26 - RenderToken
$ws.Range("E26").Value = "  +0.12%  "

# Row 27 - now LEO (was Toncoin)
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "6.06"
$ws.Range("E27").Value = "  +1.07%  "

# Row 28 - now Toncoin (was LEO)
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "3.81"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  -0.15%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.58"
$ws.Range("E31").Value = "  -6.44%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "12.79"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33 - now Hedera (was InjectiveProtocol)
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +2.95%  "

# Row 34 - now InjectiveProtocol (was Hedera)
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "46.16"
$ws.Range("E34").Value = "  +3.69%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "637.62"
$ws.Range("E35").Value = "  +2.78%  "

# Row 36 - OKB
$ws.Range("E36").Value = "  -0.84%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "0.0₃0830"
$ws.Range("E37").Value = "  -10.52%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "0.414"
$ws.Range("E38").Value = "  +2.21%  "

# Row 39 - Dai
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.14%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +2.15%  "

# Row 42 - ThetaToken
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").Value = "  -1.73%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -0.38%  "

# Row 44 - Fetch.AI
$ws.Range("D44").Value = "2.63"
$ws.Range("E44").Value = "  -1.01%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +2.70%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.879.33"
$ws.Range("E46").Value = "  +4.81%  "

# Row 47 - THORChain
$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  -3.21%  "

# Row 48 - WEMIXToken
$ws.Range("D48").Value = "2.75"
$ws.Range("E48").Value = "  -0.73%  "

# Row 49 - Monero
$ws.Range("D49").Value = "144.19"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50 - ApeXProtocol
$ws.Range("D50").Value = "3.09"
$ws.Range("E50").Value = "  -2.12%  "

# Row 51 - dogwifhat
$ws.Range("D51").Value = "2.60"
$ws.Range("E51").Value = "  -19.38%  "
